$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the "through" date references (12-08 -> 12-09)
$ws.Name = "Through 2022-12-09"

# Update the 2022 column header text
$ws.Range("I1").Value = "2022 (through 12-09)"

# Update December (row 13) and Total (row 14) figures for the 2022 column (I)
$ws.Range("I13").Value = 37
$ws.Range("I14").Value = 1553
